$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.148.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.97%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.637.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.07%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.517"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.27%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  +0.49%  "

$ws.Range("E9").Value = "  +0.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0847"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.866.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.634.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.27%  "

$ws.Range("E14").Value = "  +0.53%  "

$ws.Range("E15").Value = "  +2.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.147.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.99%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.61%  "

$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.62%  "

$ws.Range("E22").Value = "  +0.36%  "

$ws.Range("E23").Value = "  +3.60%  "

$ws.Range("E24").Value = "  -0.66%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.14%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("E27").Value = "  +2.02%  "

$ws.Range("E28").Value = "  +0.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.35%  "

$ws.Range("E30").Value = "  +0.58%  "

$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("E32").Value = "  +1.60%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.83%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.300.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.94%  "

$ws.Range("E35").Value = "  +0.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.88%  "

$ws.Range("E37").Value = "  -0.49%  "

$ws.Range("E38").Value = "  +2.02%  "

$ws.Range("E39").Value = "  +2.23%  "

$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.808"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.15%  "

$ws.Range("E42").Value = "  +5.93%  "

$ws.Range("E43").Value = "  -1.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.777.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("E45").Value = "  -0.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.71%  "

$ws.Range("E47").Value = "  -0.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0106"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0513"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.17%  "

$ws.Range("E50").Value = "  +0.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0956"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.56%  "
